$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.956.03"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.673.46"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "1.648.85"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "26.943.77"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "1.477.67"
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +7.34%  "
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.814.38"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.23%  "
